$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the RANDBETWEEN formulas in A2:A7 with fixed static values
$ws.Range("A2").Value = 1818072874.0
$ws.Range("A3").Value = 1822225502.0
$ws.Range("A4").Value = 1195849953.0
$ws.Range("A5").Value = 1206482622.0
$ws.Range("A6").Value = 1155066173.0
$ws.Range("A7").Value = 1547902419.0
